# Landscaping Data.xlsx update
# Adds 7 new rows of data (rows 387-393) for 7/4/2025 (serial date 45842),
# extends the ABS(D-E) Temp_Diff formula down through the new rows,
# and updates the active selection / window view to reflect the new data range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Prepare formatting for the new rows by copying the format of the ---
# --- last existing data row (386), which carries the date number format ---
# --- for column A and the default (General) formatting for the rest.   ---
$ws.Range("A386:T386").Copy()
$ws.Range("A387:T393").PasteSpecial(-4122)
$excel.CutCopyMode = 0

    # Row 387
    $ws.Range("A387").Value = 45842
    $ws.Range("B387").Value = "Flowering"
    $ws.Range("C387").Value = "Large"
    $ws.Range("D387").Value = 65
    $ws.Range("E387").Value = 85
    $ws.Range("F387").Formula = "=ABS(D387-E387)"
    $ws.Range("G387").Value = 0
    $ws.Range("H387").Value = 0
    $ws.Range("I387").Value = "No"
    $ws.Range("J387").Value = 2
    $ws.Range("K387").Value = "Bright"
    $ws.Range("L387").Value = 8
    $ws.Range("M387").Value = 0.48
    $ws.Range("N387").Value = 63
    $ws.Range("O387").Value = 30.13
    $ws.Range("P387").Value = 6
    $ws.Range("Q387").Value = 0.06
    $ws.Range("R387").Value = 9.9
    $ws.Range("S387").Value = 67
    $ws.Range("T387").Value = 0

    # Row 388
    $ws.Range("A388").Value = 45842
    $ws.Range("B388").Value = "Nonflowering"
    $ws.Range("C388").Value = "Medium"
    $ws.Range("D388").Value = 65
    $ws.Range("E388").Value = 85
    $ws.Range("F388").Formula = "=ABS(D388-E388)"
    $ws.Range("G388").Value = 0
    $ws.Range("H388").Value = 0
    $ws.Range("I388").Value = "No"
    $ws.Range("J388").Value = 3
    $ws.Range("K388").Value = "Bright"
    $ws.Range("L388").Value = 8
    $ws.Range("M388").Value = 0.48
    $ws.Range("N388").Value = 63
    $ws.Range("O388").Value = 30.13
    $ws.Range("P388").Value = 6
    $ws.Range("Q388").Value = 0.06
    $ws.Range("R388").Value = 9.9
    $ws.Range("S388").Value = 67
    $ws.Range("T388").Value = 0

    # Row 389
    $ws.Range("A389").Value = 45842
    $ws.Range("B389").Value = "Nonflowering"
    $ws.Range("C389").Value = "Small"
    $ws.Range("D389").Value = 65
    $ws.Range("E389").Value = 85
    $ws.Range("F389").Formula = "=ABS(D389-E389)"
    $ws.Range("G389").Value = 0
    $ws.Range("H389").Value = 0
    $ws.Range("I389").Value = "No"
    $ws.Range("J389").Value = 3
    $ws.Range("K389").Value = "Bright"
    $ws.Range("L389").Value = 8
    $ws.Range("M389").Value = 0.48
    $ws.Range("N389").Value = 63
    $ws.Range("O389").Value = 30.13
    $ws.Range("P389").Value = 6
    $ws.Range("Q389").Value = 0.06
    $ws.Range("R389").Value = 9.9
    $ws.Range("S389").Value = 67
    $ws.Range("T389").Value = 0

    # Row 390
    $ws.Range("A390").Value = 45842
    $ws.Range("B390").Value = "Nonflowering"
    $ws.Range("C390").Value = "Medium"
    $ws.Range("D390").Value = 65
    $ws.Range("E390").Value = 85
    $ws.Range("F390").Formula = "=ABS(D390-E390)"
    $ws.Range("G390").Value = 0
    $ws.Range("H390").Value = 0
    $ws.Range("I390").Value = "No"
    $ws.Range("J390").Value = 3
    $ws.Range("K390").Value = "Bright"
    $ws.Range("L390").Value = 8
    $ws.Range("M390").Value = 0.48
    $ws.Range("N390").Value = 63
    $ws.Range("O390").Value = 30.13
    $ws.Range("P390").Value = 6
    $ws.Range("Q390").Value = 0.06
    $ws.Range("R390").Value = 9.9
    $ws.Range("S390").Value = 67
    $ws.Range("T390").Value = 0

    # Row 391
    $ws.Range("A391").Value = 45842
    $ws.Range("B391").Value = "Nonflowering"
    $ws.Range("C391").Value = "Medium"
    $ws.Range("D391").Value = 65
    $ws.Range("E391").Value = 85
    $ws.Range("F391").Formula = "=ABS(D391-E391)"
    $ws.Range("G391").Value = 0
    $ws.Range("H391").Value = 0.1
    $ws.Range("I391").Value = "No"
    $ws.Range("J391").Value = 3
    $ws.Range("K391").Value = "Bright"
    $ws.Range("L391").Value = 8
    $ws.Range("M391").Value = 0.48
    $ws.Range("N391").Value = 63
    $ws.Range("O391").Value = 30.13
    $ws.Range("P391").Value = 6
    $ws.Range("Q391").Value = 0.06
    $ws.Range("R391").Value = 9.9
    $ws.Range("S391").Value = 67
    $ws.Range("T391").Value = 0

    # Row 392
    $ws.Range("A392").Value = 45842
    $ws.Range("B392").Value = "Nonflowering"
    $ws.Range("C392").Value = "Large"
    $ws.Range("D392").Value = 65
    $ws.Range("E392").Value = 85
    $ws.Range("F392").Formula = "=ABS(D392-E392)"
    $ws.Range("G392").Value = 0
    $ws.Range("H392").Value = 0.2
    $ws.Range("I392").Value = "No"
    $ws.Range("J392").Value = 4
    $ws.Range("K392").Value = "Bright"
    $ws.Range("L392").Value = 8
    $ws.Range("M392").Value = 0.48
    $ws.Range("N392").Value = 63
    $ws.Range("O392").Value = 30.13
    $ws.Range("P392").Value = 6
    $ws.Range("Q392").Value = 0.06
    $ws.Range("R392").Value = 9.9
    $ws.Range("S392").Value = 67
    $ws.Range("T392").Value = 0

    # Row 393
    $ws.Range("A393").Value = 45842
    $ws.Range("B393").Value = "Tree"
    $ws.Range("C393").Value = "Medium"
    $ws.Range("D393").Value = 65
    $ws.Range("E393").Value = 85
    $ws.Range("F393").Formula = "=ABS(D393-E393)"
    $ws.Range("G393").Value = 0
    $ws.Range("H393").Value = 0.2
    $ws.Range("I393").Value = "No"
    $ws.Range("J393").Value = 1
    $ws.Range("K393").Value = "Bright"
    $ws.Range("L393").Value = 8
    $ws.Range("M393").Value = 0.48
    $ws.Range("N393").Value = 63
    $ws.Range("O393").Value = 30.13
    $ws.Range("P393").Value = 6
    $ws.Range("Q393").Value = 0.06
    $ws.Range("R393").Value = 9.9
    $ws.Range("S393").Value = 67
    $ws.Range("T393").Value = 0


# --- The sheet's dimension updates automatically based on the used range. ---

# --- Scroll the view and set the active selection to match where the ---
# --- author was last working (mirrors the sheetView/selection change). ---
$win = $excel.ActiveWindow
$win.ScrollRow = 383
$win.ScrollColumn = 1
[void]$ws.Range("Q387:Q393").Select()

# --- Reposition the workbook window, matching the saved window state. ---
$win.Left = 804
$win.Top = 3360
